$d = $word.ActiveDocument

$replacements = @(
    @("444÷3=", "295÷9="),
    @("621÷8=", "354÷9="),
    @("842÷8=", "448÷7="),
    @("917÷6=", "765÷7="),
    @("115÷9=", "998÷5="),
    @("556÷3=", "811÷7="),
    @("439÷8=", "950÷6="),
    @("906÷9=", "825÷6="),
    @("338÷9=", "921÷5="),
    @("797÷4=", "925÷7="),
    @("622÷4=", "753÷5="),
    @("665÷3=", "739÷8="),
    @("744÷5=", "337÷2="),
    @("520÷8=", "654÷5="),
    @("460÷7=", "236÷5="),
    @("382÷2=", "132÷7="),
    @("193÷6=", "178÷5="),
    @("700÷9=", "473÷3="),
    @("398÷6=", "400÷7="),
    @("542÷8=", "775÷4="),
    @("190÷3=", "578÷3="),
    @("627÷4=", "323÷5="),
    @("667÷3=", "327÷5="),
    @("375÷9=", "533÷7="),
    @("528÷5=", "356÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
